$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the number format on S12 (it was using the wrong "date only" style;
# align it with the rest of the Date column which uses date+time).
$ws.Range("S12").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new data row (row 13) with the latest bunker prices.
$rowValues = @(518,570,535,520,638,525,640,658,580,525,536,509,590,515,637,762,563,655,45734,580,552,595,498,509,545,752,536,560,520,639,591.5,564,530,564,881,637,505,622,544,511,529,509,498,495,525,550,495,534)

for ($i = 0; $i -lt $rowValues.Length; $i++) {
    $ws.Cells.Item(13, $i + 1).Value = $rowValues[$i]
}

# New row's Date cell keeps the "date only" number format (matches the diff).
$ws.Range("S13").NumberFormat = "YYYY-MM-DD"
